$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows: origin/destination swaps (JFK/ORY retired, AMS/LDN introduced) ---
$ws.Cells.Item(4, 3).Value = "AMS"    # C4: JFK -> AMS

$ws.Cells.Item(5, 3).Value = "AMS"    # C5: CDG -> AMS
$ws.Cells.Item(5, 4).Value = "CDG"    # D5: JFK -> CDG

$ws.Cells.Item(6, 4).Value = "AMS"    # D6: VCE -> AMS

$ws.Cells.Item(8, 3).Value = "CDG"    # C8: VCE -> CDG
$ws.Cells.Item(8, 4).Value = "AMS"    # D8: ORY -> AMS

$ws.Cells.Item(9, 3).Value = "AMS"    # C9: AMS -> AMS (stays AMS, different shared-string slot)
$ws.Cells.Item(9, 4).Value = "LDN"    # D9: JFK -> LDN

# --- Append new demo flight / booking rows (10-15) ---

# Row 10
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "ON TIME"
$ws.Cells.Item(10, 3).Value = "AMS"
$ws.Cells.Item(10, 4).Value = "LDN"
$ws.Cells.Item(10, 5).Value = 100
$ws.Cells.Item(10, 6).Value = 45280
$ws.Cells.Item(10, 7).Value = 0.41666666666666669
$ws.Cells.Item(10, 7).NumberFormat = "h:mm"

# Row 11
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "ON TIME"
$ws.Cells.Item(11, 3).Value = "AMS"
$ws.Cells.Item(11, 4).Value = "LDN"
$ws.Cells.Item(11, 5).Value = 200
$ws.Cells.Item(11, 6).Value = 45281
$ws.Cells.Item(11, 7).Value = 0.4375
$ws.Cells.Item(11, 7).NumberFormat = "h:mm"

# Row 12
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "DELAYED"
$ws.Cells.Item(12, 3).Value = "LDN"
$ws.Cells.Item(12, 4).Value = "AMS"
$ws.Cells.Item(12, 5).Value = 150
$ws.Cells.Item(12, 6).Value = 45285
$ws.Cells.Item(12, 7).Value = 0.95833333333333337
$ws.Cells.Item(12, 7).NumberFormat = "h:mm"

# Row 13
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "ON TIME"
$ws.Cells.Item(13, 3).Value = "LDN"
$ws.Cells.Item(13, 4).Value = "AMS"
$ws.Cells.Item(13, 5).Value = 200
$ws.Cells.Item(13, 6).Value = 45287
$ws.Cells.Item(13, 7).Value = 0.375
$ws.Cells.Item(13, 7).NumberFormat = "h:mm"

# Row 14
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = "ON TIME"
$ws.Cells.Item(14, 3).Value = "AMS"
$ws.Cells.Item(14, 4).Value = "LDN"
$ws.Cells.Item(14, 5).Value = 250
$ws.Cells.Item(14, 6).Value = 45290
$ws.Cells.Item(14, 7).Value = 0.5
$ws.Cells.Item(14, 7).NumberFormat = "h:mm"

# Row 15
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "ON TIME "
$ws.Cells.Item(15, 3).Value = "AMS"
$ws.Cells.Item(15, 4).Value = "CDG"
$ws.Cells.Item(15, 5).Value = 150
$ws.Cells.Item(15, 6).Value = 45290
$ws.Cells.Item(15, 7).Value = 0.58333333333333337
$ws.Cells.Item(15, 7).NumberFormat = "h:mm"

# --- Update selection to reflect where the user ended up after entering data ---
$ws.Range("A16").Select()
